# Insert a new row for the "DeepCNN" model above the existing "GBDT" row
# (originally row 15), shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = "DeepCNN"
